$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto snapshot values.
# NumberFormat is forced to Text ("@") before assignment so numeric-looking strings
# (including ones with significant trailing zeros, e.g. "37.20") are preserved verbatim
# instead of being coerced into floating point numbers.
$updates = @{
    'D2' = '308.09'
    'E2' = '-0.46%'
    'D3' = '37.20'
    'D4' = '5.124'
    'E4' = '-0.10%'
    'D5' = '0.07834'
    'E5' = '0.91%'
    'D6' = '8.251'
    'E6' = '0.62%'
    'D7' = '1.890'
    'E7' = '0.68%'
    'D8' = '2.943'
    'E8' = '-3.54%'
    'D9' = '0.9210'
    'E9' = '0.35%'
    'D10' = '0.1082'
    'E10' = '-9.02%'
    'D11' = '0.1896'
    'E11' = '0.20%'
    'D12' = '0.08852'
    'E12' = '-8.44%'
    'D13' = '0.03340'
    'E13' = '-2.85%'
    'D14' = '0.09591'
    'E14' = '-1.04%'
    'D15' = '0.001375'
    'E15' = '0.18%'
    'D16' = '0.005737'
    'E16' = '-0.77%'
    'D17' = '3.412'
    'E17' = '-3.56%'
    'D18' = '4.394'
    'E18' = '-0.16%'
    'D19' = '0.3425'
    'E19' = '0.77%'
    'D20' = '6.294'
    'E20' = '19.87%'
    'E21' = '1.46%'
    'D22' = '0.2415'
    'E22' = '-6.88%'
    'D23' = '0.04359'
    'E23' = '0.70%'
    'D24' = '0.001193'
    'E24' = '-0.55%'
    'D25' = '0.004270'
    'E25' = '0.44%'
    'E26' = '7.49%'
    'E27' = '-98.11%'
    'D39' = '0.02175'
    'E39' = '5.47%'
    'D40' = '0.05037'
    'E40' = '0.38%'
    'D41' = '0.007573'
    'E41' = '-1.25%'
    'D42' = '0.1351'
    'E42' = '0.75%'
    'D43' = '0.008664'
    'D44' = '0.001979'
    'E44' = '-8.90%'
    'D45' = '0.007899'
    'E45' = '-9.77%'
    'D46' = '0.00006527'
    'E46' = '-2.73%'
    'E47' = '-0.18%'
    'D48' = '0.003295'
    'E48' = '12.12%'
    'E49' = '-16.58%'
    'D50' = '0.00002099'
    'E50' = '-0.18%'
    'D51' = '0.0001999'
    'E51' = '-0.18%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"
